# Lista_Exemplo.xlsx update
# - Disciplinas: "Codigo" (col C) renumbered from 1-7 to text "8"-"14"
# - Turmas: "Nome da Turma" (col B) renumbered from T001-T012 to T013-T024
# - Usuarios: "Matricula" (col B) updated to new enrollment numbers (now stored as text)
# - Vinculos: "Matricula" (col B) updated to match the new Usuarios numbers (stays numeric)
# Plus the small UI/view bookkeeping (selections + active sheet) that Excel
# records when a user edits these sheets in turn, finishing on "Vinculos".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Turmas — Nome da Turma T001..T012 -> T013..T024
# ---------------------------------------------------------------------------
$wsTurmas = $wb.Worksheets.Item("Turmas")

$wsTurmas.Range("B2").Value = "T013"
$wsTurmas.Range("B3").Value = "T014"
$wsTurmas.Range("B4").Value = "T015"
$wsTurmas.Range("B5").Value = "T016"
$wsTurmas.Range("B6").Value = "T017"
$wsTurmas.Range("B7").Value = "T018"
$wsTurmas.Range("B8").Value = "T019"
$wsTurmas.Range("B9").Value = "T020"
$wsTurmas.Range("B10").Value = "T021"
$wsTurmas.Range("B11").Value = "T022"
$wsTurmas.Range("B12").Value = "T023"
$wsTurmas.Range("B13").Value = "T024"

$wsTurmas.Range("B16").Select()

# ---------------------------------------------------------------------------
# 2) Disciplinas — Codigo column becomes text "8".."14"
# ---------------------------------------------------------------------------
$wsDisciplinas = $wb.Worksheets.Item("Disciplinas")

$wsDisciplinas.Range("C2").Value = "8"
$wsDisciplinas.Range("C3").Value = "9"
$wsDisciplinas.Range("C4").Value = "10"
$wsDisciplinas.Range("C5").Value = "11"
$wsDisciplinas.Range("C6").Value = "12"
$wsDisciplinas.Range("C7").Value = "13"
$wsDisciplinas.Range("C8").Value = "14"

# old leftover formatted cell at C18 goes away, a new one shows up at D12
$wsDisciplinas.Range("C18").Clear()
$wsDisciplinas.Range("D12").Font.Underline = $true

$wsDisciplinas.Range("D12").Select()

# ---------------------------------------------------------------------------
# 3) Usuarios — Matricula renumbered (now stored as text)
# ---------------------------------------------------------------------------
$wsUsuarios = $wb.Worksheets.Item("Usuarios")

$wsUsuarios.Range("B2").Value = "20251853"
$wsUsuarios.Range("B3").Value = "20255832"
$wsUsuarios.Range("B4").Value = "20249245"
$wsUsuarios.Range("B5").Value = "20254321"
$wsUsuarios.Range("B6").Value = "20236123"
$wsUsuarios.Range("B7").Value = "20232935"

# a new formatted (but empty) row shows up below the data
$wsUsuarios.Range("C11").Font.Underline = $true

$wsUsuarios.Range("B7").Select()

# ---------------------------------------------------------------------------
# 4) Vinculos — Matricula updated to match Usuarios (stays numeric)
# ---------------------------------------------------------------------------
$wsVinculos = $wb.Worksheets.Item("Vinculos")

$wsVinculos.Range("B2").Value = 20251853
$wsVinculos.Range("B3").Value = 20236123
$wsVinculos.Range("B4").Value = 20232935

$wsVinculos.Range("B4").Select()

# Vinculos ends up the active sheet/tab after these edits
$wsVinculos.Activate()
